$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 457, shifting existing rows 457:485 down to 458:486.
$ws.Rows.Item(457).Insert()

# Populate the newly inserted row 457 with the new weekly record.
$ws.Cells.Item(457, 1).Value = 10
$ws.Cells.Item(457, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(457, 3).Value = "La Araucanía"
$ws.Cells.Item(457, 4).Value = 45021
$ws.Cells.Item(457, 5).Value = 9
$ws.Cells.Item(457, 6).Value = 100112009
$ws.Cells.Item(457, 7).Value = "Acelga"
$ws.Cells.Item(457, 8).Value = "Sin especificar"
$ws.Cells.Item(457, 9).Value = "Primera"
$ws.Cells.Item(457, 10).Value = 65
$ws.Cells.Item(457, 11).Value = 8000
$ws.Cells.Item(457, 12).Value = 8000
$ws.Cells.Item(457, 13).Value = 8000
$ws.Cells.Item(457, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(457, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(457, 16).Value = 667
$ws.Cells.Item(457, 17).Value = 12
$ws.Cells.Item(457, 18).Value = "Hortaliza"
